# table ss with pre-experiment data
#
# The "exit" and "basal" sheets hold raw balance-test statistics that get
# refreshed with pre-experiment data. The "balance_response" sheet pulls
# those numbers in through formulas (=ROUND(basal!..,2) etc.) and, in the
# source commit, those formula caches were NOT refreshed when the raw data
# was pasted in (fullCalcOnLoad was set instead, so Excel recalculates the
# next time the file is opened). We reproduce that by switching to manual
# calculation before writing the raw values, so the dependent formulas on
# balance_response keep their old cached results.
$excel.Calculation = -4135  # xlCalculationManual

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "exit" sheet (sheet1) - new pre-experiment values for B2:E8
# ---------------------------------------------------------------------
$exit = $wb.Worksheets.Item("exit")

$exit.Cells.Item(2,2).Value = 2157.4058537363426
$exit.Cells.Item(2,3).Value = 2367.6243093922653
$exit.Cells.Item(2,4).Value = 2171.5569770901516
$exit.Cells.Item(2,5).Value = 0.17228973589263194

$exit.Cells.Item(3,2).Value = 30.915340842769982
$exit.Cells.Item(3,3).Value = 154.83134488635548
$exit.Cells.Item(3,4).Value = 31.807323205390613

$exit.Cells.Item(4,2).Value = 0.18231118909003907
$exit.Cells.Item(4,3).Value = 0.19337016574585636
$exit.Cells.Item(4,4).Value = 0.18305563820291579
$exit.Cells.Item(4,5).Value = 0.71034364024861918

$exit.Cells.Item(5,2).Value = 0.024499682973906006
$exit.Cells.Item(5,3).Value = 0.036351262642859185
$exit.Cells.Item(5,4).Value = 0.024353374831239075

$exit.Cells.Item(6,2).Value = 33.340425531914896
$exit.Cells.Item(6,3).Value = 29.29032258064516
$exit.Cells.Item(6,4).Value = 33.031941031941031
$exit.Cells.Item(6,5).Value = 0.24593205180445749

$exit.Cells.Item(7,2).Value = 1.0865666714159887
$exit.Cells.Item(7,3).Value = 3.3579286845159224
$exit.Cells.Item(7,4).Value = 1.0362073711473245

$exit.Cells.Item(8,2).Value = 12539
$exit.Cells.Item(8,3).Value = 905
$exit.Cells.Item(8,4).Value = 13444

# ---------------------------------------------------------------------
# "basal" sheet (sheet2) - new pre-experiment values for B2:E8
# ---------------------------------------------------------------------
$basal = $wb.Worksheets.Item("basal")

$basal.Cells.Item(2,2).Value = 2262.1439973395409
$basal.Cells.Item(2,3).Value = 2145.4579860113058
$basal.Cells.Item(2,4).Value = 2171.5569770901516
$basal.Cells.Item(2,5).Value = 0.077363019270749664

$basal.Cells.Item(3,2).Value = 59.211959578633127
$basal.Cells.Item(3,3).Value = 35.445186055768197
$basal.Cells.Item(3,4).Value = 31.80732320539061

$basal.Cells.Item(4,2).Value = 0.18190887928167609
$basal.Cells.Item(4,3).Value = 0.18338603046852545
$basal.Cells.Item(4,4).Value = 0.18305563820291579
$basal.Cells.Item(4,5).Value = 0.94953992198672665

$basal.Cells.Item(5,2).Value = 0.031453663923441053
$basal.Cells.Item(5,3).Value = 0.024509689361658406
$basal.Cells.Item(5,4).Value = 0.024353374831239075

$basal.Cells.Item(6,2).Value = 32.127450980392155
$basal.Cells.Item(6,3).Value = 33.334426229508196
$basal.Cells.Item(6,4).Value = 33.031941031941031
$basal.Cells.Item(6,5).Value = 0.61423991149630908

$basal.Cells.Item(7,2).Value = 2.0762812455551489
$basal.Cells.Item(7,3).Value = 1.1972335117199615
$basal.Cells.Item(7,4).Value = 1.0362073711473245

$basal.Cells.Item(8,2).Value = 3007
$basal.Cells.Item(8,3).Value = 10437
$basal.Cells.Item(8,4).Value = 13444

# ---------------------------------------------------------------------
# "balance_response" sheet (sheet3) - relabel the pawns row, resize the
# label column to fit the new text, and move the active selection back
# to A2 (matching the saved view state).
# ---------------------------------------------------------------------
$balance = $wb.Worksheets.Item("balance_response")
$balance.Range("A8").Value = "Number of branch-day pawns"
$balance.Columns.Item(1).AutoFit()

$balance.Activate()
$balance.Range("A2:H10").Select()
